# Apply "User data with filename fix" changes:
#  - Update "Units of measure" note on the Notes sheet
#  - Populate the Data sheet with country/year/value rows

$wb = $excel.ActiveWorkbook
$notes = $wb.Worksheets.Item("Notes")
$data = $wb.Worksheets.Item("Data")

# Fix the units-of-measure note text (was "NA")
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"

# Fill in the Data sheet rows
$rows = @(
    @("IQ", "Iraq", 2014, 11944.45),
    @("JO", "Jordan", 2015, 55000),
    @("MG", "Madagascar", 2015, 6200),
    @("MW", "Malawi", 2015, 5500),
    @("middle-east", "Middle East, regional", 2014, 7249.05),
    @("PS", "Palestine", 2014, 65900.43),
    @("PS", "Palestine", 2015, 490000),
    @("PH", "Philippines", 2014, 21829.52),
    @("RS", "Serbia", 2014, 21911.89),
    @("SY", "Syria", 2014, 20593.88),
    @("SY", "Syria", 2015, 11000),
    @("UG", "Uganda", 2014, 4942.53)
)

$r = 2
foreach ($row in $rows) {
    $data.Cells.Item($r, 1).Value = $row[0]
    $data.Cells.Item($r, 2).Value = $row[1]
    $data.Cells.Item($r, 3).Value = $row[2]
    $data.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
